$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "2025/12/03 14:00"
$ws.Range("B22").Value = "-"
$ws.Range("C22").Value = "-"
$ws.Range("D22").Value = "-"
$ws.Range("E22").Value = "-"
$ws.Range("F22").Value = "-"
$ws.Range("G22").Value = "-"
